$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

$ws.Range("A6").Value = "Bengali"
$ws.Range("B6").Value = "Sanjib Roy"
$ws.Range("C6").Value = "sanjibroy0098@gmail.com"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
